$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Card18")

# --- Row 16: the previously-blank log columns get the literal placeholder
#     text "nan" (matching the convention already used elsewhere in this
#     sheet for service-log rows with missing readings).
$ws.Range("B16").Value = "nan"
$ws.Range("C16").Value = "nan"
$ws.Range("D16").Value = "nan"
$ws.Range("E16").Value = "nan"
$ws.Range("F16").Value = "nan"
$ws.Range("G16").Value = "nan"
$ws.Range("H16").Value = "nan"
$ws.Range("I16").Value = "nan"
$ws.Range("J16").Value = "nan"
$ws.Range("K16").Value = "nan"
$ws.Range("M16").Value = "nan"
$ws.Range("N16").Value = "nan"

# --- Row 17: a brand-new service-log entry for Card18.
# A17 and M17 look numeric ("18" / "590.1") but the sheet stores every
# column as text, so force a text number format before assigning them
# to avoid Excel auto-converting the literal into a real number.
$ws.Range("A17").NumberFormat = "@"
$ws.Range("A17").Value = "18"

$ws.Range("L17").Value = "29\4\2025"

$ws.Range("M17").NumberFormat = "@"
$ws.Range("M17").Value = "590.1"

$ws.Range("O17").Value = "تم تغيير الجرائد الاماميه (1_2_4__5_7_8) ومعيارته وسن السليندر"
$ws.Range("P17").Value = "الخبير"
